# Update gh-pages output data (苏州-漫展信息) to the values generated at 456a3b4.
#
# Sheet "展览" (rId1) and sheet "全部类型" (rId4) both list the same set of
# events (the latter is a superset covering every category), so each event's
# "想去人数" (F column) count and, for one event, its cover-image URL (I
# column) need to be bumped identically in both sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- "展览" sheet: row -> new F value -------------------------------------
$exhibitCounts = @{
    2  = 1177
    4  = 275
    7  = 12304
    8  = 60
    11 = 149
    12 = 12106
    13 = 4816
    14 = 4676
    21 = 361
    23 = 72
}
foreach ($row in $exhibitCounts.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitCounts[$row]
}
$wsExhibit.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202408/f77bW6VT1724292729739.jpeg"

# --- "全部类型" sheet: row -> new F value -----------------------------------
$allCounts = @{
    2  = 1177
    4  = 275
    9  = 12304
    10 = 60
    13 = 149
    14 = 12106
    15 = 4816
    16 = 4676
    23 = 361
    25 = 72
}
foreach ($row in $allCounts.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allCounts[$row]
}
$wsAll.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202408/f77bW6VT1724292729739.jpeg"
